$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

# Update the list of predictor variables (x) - dictionary was regenerated
$ws.Range("B2").Value = "Ones, NEMScr, LangScr, MathScr, SchoolRegion_2, MotherEd_4, MotherEd_7, PostulationType_1"

# Widen column B so the longer text still fits (bestFit-style width)
$ws.Columns.Item(2).ColumnWidth = 92

# Update the recalculated model metrics (pyplot/model state fix changed the numbers)
$ws.Range("B4").Value = 0.69444444444444442
$ws.Range("B5").Value = 0.98648648648648651
$ws.Range("B6").Value = 0.69523809523809521
$ws.Range("B7").Value = 0.81564245810055858
